# Apply updated crypto price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''59.112.52'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '''2.508.20'
$ws.Range('E3').Value = '  +2.05%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''541.59'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').Value = '''143.75'
$ws.Range('E6').Value = '  -3.06%  '
$ws.Range('D7').Value = '''0.996'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').Value = '''2.534.55'
$ws.Range('E9').Value = '  +2.30%  '
$ws.Range('E10').Value = '  +1.06%  '
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('E12').Value = '  +3.45%  '
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '''2.952.42'
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('D15').Value = '''23.57'
$ws.Range('E15').Value = '  -2.60%  '
$ws.Range('D16').Value = '''59.049.24'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D18').Value = '''2.527.82'
$ws.Range('E18').Value = '  +0.35%  '
$ws.Range('D19').Value = '''11.21'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '''4.29'
$ws.Range('E20').Value = '  -1.84%  '
$ws.Range('D21').Value = '''324.53'
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('E22').Value = '  +2.96%  '
$ws.Range('D23').Value = '''5.78'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('E25').Value = '  -5.35%  '
$ws.Range('E26').Value = '  +1.30%  '
$ws.Range('D27').Value = '''0.996'
$ws.Range('E27').Value = '  +1.61%  '
$ws.Range('D28').Value = '''7.89'
$ws.Range('E28').Value = '  +1.76%  '
$ws.Range('D29').Value = '''0.0₃0781'
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').Value = '''6.63'
$ws.Range('E31').Value = '  -1.65%  '
$ws.Range('D32').Value = '''1.20'
$ws.Range('E32').Value = '  -8.44%  '
$ws.Range('D33').Value = '''0.997'
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('E34').Value = '  +6.11%  '
$ws.Range('D35').Value = '''158.44'
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('D36').Value = '''18.67'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').Value = '''4.37'
$ws.Range('E37').Value = '  -4.14%  '
$ws.Range('D38').Value = '''1.62'
$ws.Range('E38').Value = '  -7.07%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''5.66'
$ws.Range('E39').Value = '  -4.57%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '''36.94'
$ws.Range('E40').Value = '  +0.95%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = '''296.72'
$ws.Range('E41').Value = '  -7.35%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '''3.70'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('D44').Value = '''0.994'
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('D45').Value = '''0.601'
$ws.Range('E45').Value = '  +2.54%  '
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('D48').Value = '''18.69'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('D49').Value = '''122.63'
$ws.Range('E49').Value = '  +0.45%  '
$ws.Range('D50').Value = '''0.0516'
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('E51').Value = '  -1.04%  '
